$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late" / variable instalment column),
# shifting the existing "Late" and "Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()

# Match the width Excel applies to the newly inserted column (renders as width="10").
$ws.Columns("N:N").ColumnWidth = 9.140625

# Restore the active selection as recorded after the edit.
$ws.Range("S9").Select() | Out-Null
